$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 10417083
$ws.Range("I2").Value = 16666725
$ws.Range("J2").Value = 1013.3333
$ws.Range("K2").Value = 16666725
$ws.Range("L2").Value = 1013.3333
$ws.Range("M2").Value = -16666612
$ws.Range("N2").Value = -1239.3333
# Row 17
$ws.Range("H17").Value = 36659.535
$ws.Range("J17").Value = 36659.535
$ws.Range("L17").Value = 109978.605
$ws.Range("N17").Value = -110314.605
# Row 18
$ws.Range("H18").Value = 816.6667
$ws.Range("I18").Value = 800
$ws.Range("J18").Value = 850
$ws.Range("K18").Value = 800
$ws.Range("L18").Value = 850
$ws.Range("M18").Value = -516
$ws.Range("N18").Value = -1418
# Row 40
$ws.Range("H40").Value = 53159850
$ws.Range("I40").Value = 1469.0769
$ws.Range("J40").Value = 168336340
$ws.Range("K40").Value = 1469.0769
$ws.Range("L40").Value = 168336340
$ws.Range("M40").Value = -1294.0769
$ws.Range("N40").Value = -168336690
# Row 41
$ws.Range("H41").Value = 144.27272
$ws.Range("I41").Value = 138.7
$ws.Range("J41").Value = 148.91667
$ws.Range("K41").Value = 138.7
$ws.Range("L41").Value = 148.91667
$ws.Range("M41").Value = 301.3
$ws.Range("N41").Value = -1028.91667
# Row 129
$ws.Range("H129").Value = 903.63416
$ws.Range("J129").Value = 950
$ws.Range("L129").Value = 2850
$ws.Range("N129").Value = -12850

$ws = $wb.Worksheets.Item("ARM")
# Row 26
$ws.Range("H26").Value = 2160
$ws.Range("I26").Value = 840
$ws.Range("J26").Value = 4800
$ws.Range("K26").Value = 840
$ws.Range("L26").Value = 4800
$ws.Range("M26").Value = -510
$ws.Range("N26").Value = -5460
# Row 32
$ws.Range("H32").Value = 1486.54
$ws.Range("I32").Value = 1453.9656
$ws.Range("J32").Value = 1704.5385
$ws.Range("K32").Value = 1453.9656
$ws.Range("L32").Value = 1704.5385
$ws.Range("M32").Value = -1166.9656
$ws.Range("N32").Value = -2278.5385
# Row 74
$ws.Range("H74").Value = 838
$ws.Range("I74").Value = 847.6486
$ws.Range("K74").Value = 847.6486
$ws.Range("M74").Value = 26.35140000000001
# Row 77
$ws.Range("H77").Value = 838
$ws.Range("I77").Value = 847.6486
$ws.Range("K77").Value = 4238.243
$ws.Range("M77").Value = 129.7569999999996

$ws = $wb.Worksheets.Item("BSM")
# Row 33
$ws.Range("H33").Value = 27100
$ws.Range("I33").Value = 20833.334
$ws.Range("J33").Value = 30233.334
$ws.Range("K33").Value = 20833.334
$ws.Range("L33").Value = 30233.334
$ws.Range("M33").Value = -20497.334
$ws.Range("N33").Value = -30905.334
# Row 86
$ws.Range("H86").Value = 970714
$ws.Range("I86").Value = 1932.4
$ws.Range("J86").Value = 1662700.9
$ws.Range("K86").Value = 1932.4
$ws.Range("L86").Value = 1662700.9
$ws.Range("M86").Value = -809.4000000000001
$ws.Range("N86").Value = -1664946.9
# Row 89
$ws.Range("H89").Value = 970714
$ws.Range("I89").Value = 1932.4
$ws.Range("J89").Value = 1662700.9
$ws.Range("K89").Value = 9662
$ws.Range("L89").Value = 8313504.5
$ws.Range("M89").Value = -4046
$ws.Range("N89").Value = -8324736.5

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 1267.6666
$ws.Range("I31").Value = 1005.70215
$ws.Range("J31").Value = 1664.8387
$ws.Range("K31").Value = 1005.70215
$ws.Range("L31").Value = 1664.8387
$ws.Range("M31").Value = -710.70215
$ws.Range("N31").Value = -2254.8387
# Row 32
$ws.Range("H32").Value = 6127.5
$ws.Range("I32").Value = 6127.5
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 6127.5
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -5811.5
$ws.Range("N32").ClearContents()
# Row 34
$ws.Range("H34").Value = 1267.6666
$ws.Range("I34").Value = 1005.70215
$ws.Range("J34").Value = 1664.8387
$ws.Range("K34").Value = 1005.70215
$ws.Range("L34").Value = 1664.8387
$ws.Range("M34").Value = -803.70215
$ws.Range("N34").Value = -2068.8387
# Row 35
$ws.Range("H35").Value = 11697.917
$ws.Range("I35").Value = 11697.917
$ws.Range("K35").Value = 11697.917
$ws.Range("M35").Value = -11403.917
# Row 38
$ws.Range("H38").Value = 19199.666
# Row 41
$ws.Range("H41").Value = 26677
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 26677
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 26677
$ws.Range("M41").ClearContents()
$ws.Range("N41").Value = -27533
# Row 46
$ws.Range("H46").Value = 19199.666
# Row 50
$ws.Range("H50").Value = 45000
$ws.Range("J50").Value = 45000
$ws.Range("L50").Value = 45000
$ws.Range("N50").Value = -46250
# Row 51
$ws.Range("H51").Value = 19332.666
$ws.Range("J51").Value = 24099
$ws.Range("L51").Value = 24099
$ws.Range("N51").Value = -25571
# Row 60
$ws.Range("H60").Value = 15799
$ws.Range("I60").Value = 12031
$ws.Range("J60").Value = 27103
$ws.Range("K60").Value = 12031
$ws.Range("L60").Value = 27103
$ws.Range("M60").Value = -11520
$ws.Range("N60").Value = -28125
# Row 61
$ws.Range("H61").Value = 19332.666
$ws.Range("J61").Value = 24099
$ws.Range("L61").Value = 24099
$ws.Range("N61").Value = -24795
# Row 74
$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()
# Row 77
$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()
# Row 94
$ws.Range("H94").Value = 7831.0713
$ws.Range("I94").Value = 20380.4
$ws.Range("J94").Value = 859.2222
$ws.Range("K94").Value = 20380.4
$ws.Range("L94").Value = 859.2222
$ws.Range("M94").Value = -19929.4
$ws.Range("N94").Value = -1761.2222
# Row 107
$ws.Range("H107").Value = 856
$ws.Range("I107").Value = 725.8182
$ws.Range("J107").Value = 1333.3334
$ws.Range("K107").Value = 725.8182
$ws.Range("L107").Value = 1333.3334
$ws.Range("M107").Value = 1194.1818
$ws.Range("N107").Value = -5173.3334

$ws = $wb.Worksheets.Item("CUL")
# Row 131
$ws.Range("H131").Value = 14587276
$ws.Range("J131").Value = 7467164
$ws.Range("L131").Value = 22401492
$ws.Range("N131").Value = -22411572

$ws = $wb.Worksheets.Item("GSM")
# Row 132
$ws.Range("H132").Value = 16242.286
$ws.Range("I132").Value = 20759.6
$ws.Range("J132").Value = 13732.667
$ws.Range("K132").Value = 62278.8
$ws.Range("L132").Value = 41198.001
$ws.Range("M132").Value = -59748.8
$ws.Range("N132").Value = -46258.001

$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 1195.3529
$ws.Range("I46").Value = 614.8
$ws.Range("J46").Value = 2024.7142
$ws.Range("K46").Value = 2024.7142
$ws.Range("L46").Value = 2024.7142
$ws.Range("M46").Value = -426.8
$ws.Range("N46").Value = -2400.7142
